$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.330.83'
$ws.Range("E2").Value = '  -0.10%  '
$ws.Range("D3").Value = '1.932.48'
$ws.Range("E3").Value = '  -0.28%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7470'
$ws.Range("E5").Value = '  +3.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '248.59'
$ws.Range("E6").Value = '  -0.74%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '28.28'
$ws.Range("E8").Value = '  -0.28%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3203'
$ws.Range("E9").Value = '  -3.89%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07096'
$ws.Range("E10").Value = '  -1.85%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7878'
$ws.Range("E11").Value = '  -2.94%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08000'
$ws.Range("D13").Value = '1.930.68'
$ws.Range("E13").Value = '  -0.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.379'
$ws.Range("E14").Value = '  -1.75%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '94.54'
$ws.Range("E15").Value = '  +0.30%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.63'
$ws.Range("E16").Value = '  -2.35%  '
$ws.Range("D17").Value = '30.333.82'
$ws.Range("E17").Value = '  -0.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '253.16'
$ws.Range("E18").Value = '  +1.72%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008036'
$ws.Range("E19").Value = '  -2.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.807'
$ws.Range("E20").Value = '  -1.73%  '
$ws.Range("D21").Value = '2.186.59'
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("E23").Value = '  -0.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.810'
$ws.Range("E24").Value = '  -2.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.579'
$ws.Range("E25").Value = '  -1.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.48'
$ws.Range("E26").Value = '  +0.72%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.336'
$ws.Range("E27").Value = '  -2.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.10'
$ws.Range("E28").Value = '  -0.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1331'
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.360'
$ws.Range("E30").Value = '  +1.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.530'
$ws.Range("E31").Value = '  -2.67%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.444'
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.147'
$ws.Range("E33").Value = '  -1.14%  '
$ws.Range("E34").Value = '  -1.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.283'
$ws.Range("E35").Value = '  -0.43%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7501'
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.769'
$ws.Range("E37").Value = '  +0.67%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01970'
$ws.Range("E38").Value = '  -0.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.803'
$ws.Range("E39").Value = '  -1.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '78.21'
$ws.Range("E40").Value = '  -3.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.416'
$ws.Range("E41").Value = '  -0.62%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4510'
$ws.Range("E42").Value = '  -0.76%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.993'
$ws.Range("E43").Value = '  -2.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8433'
$ws.Range("E44").Value = '  -0.60%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.000'
$ws.Range("E45").Value = '  -0.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.47'
$ws.Range("E46").Value = '  +0.34%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.831'
$ws.Range("E47").Value = '  +0.09%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.548'
$ws.Range("E48").Value = '  +1.45%  '
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '988.17'
$ws.Range("E49").Value = '  +12.91%  '
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.49'
$ws.Range("E50").Value = '  +1.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1206'
$ws.Range("E51").Value = '  +5.93%  '
